# Bump "Version 1." to "Version 2." reproducing the exact run layout
# produced by Word when the edit is made in place (spell-check run split
# on "Version" -> "Versi"/"on", and the trailing "." typed back in after
# the "_GoBack" bookmark once "1" was replaced by "2").

$d = $word.ActiveDocument

# --- Split "Version" into "Versi" + "on" -------------------------------
# Delete "on" (chars 5-6) and retype it; then mint+drop a bookmark right
# at the seam so the two runs aren't coalesced back into one "Version".
$rVersion = $d.Range(5, 7)
$rVersion.Text = ""
$rOn = $d.Range(5, 5)
$rOn.InsertAfter("on")

$seam1 = $d.Range(5, 5)
$d.Bookmarks.Add("_tmp_seam1", $seam1) | Out-Null
$d.Bookmarks("_tmp_seam1").Delete()

# --- "1" -> "2" ----------------------------------------------------------
$rDigit = $d.Range(8, 9)
$rDigit.Text = "2"

# --- Split the trailing "." into its own run after the bookmark --------
$rPeriod = $d.Range(9, 10)
$rPeriod.Text = ""
$rNewPeriod = $d.Range(9, 9)
$rNewPeriod.InsertAfter(".")
